# chore: update Sheets via scheduled runner
#
# Refreshes the cached market-price derived columns (H:N -
# currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
# LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ) on the
# per-class Leve-profit sheets (ALC, ARM, CRP, CUL, GSM, LTW, WVR) with
# newly pulled data. Columns A:G (leve metadata) are untouched.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 6129.8335
$ws.Range("I64").Value = 4537.5
$ws.Range("J64").Value = 6926
$ws.Range("K64").Value = 4537.5
$ws.Range("L64").Value = 6926
$ws.Range("M64").Value = -4289.5
$ws.Range("N64").Value = -7422

$ws.Range("H67").Value = 6129.8335
$ws.Range("I67").Value = 4537.5
$ws.Range("J67").Value = 6926
$ws.Range("K67").Value = 4537.5
$ws.Range("L67").Value = 6926
$ws.Range("M67").Value = -3679.5
$ws.Range("N67").Value = -8642

$ws.Range("H74").Value = 4436.625
$ws.Range("I74").Value = 3698.6
$ws.Range("K74").Value = 3698.6
$ws.Range("M74").Value = -2762.6

$ws.Range("H76").Value = 2500
$ws.Range("I76").Value = 2500
$ws.Range("K76").Value = 2500
$ws.Range("M76").Value = -2185

$ws.Range("H77").Value = 4436.625
$ws.Range("I77").Value = 3698.6
$ws.Range("K77").Value = 18493
$ws.Range("M77").Value = -13813

$ws.Range("H79").Value = 2500
$ws.Range("I79").Value = 2500
$ws.Range("K79").Value = 2500
$ws.Range("M79").Value = -1408

$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()

$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()

$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()

$ws.Range("H141").Value = 63116.934
$ws.Range("I141").Value = 72211.92
$ws.Range("K141").Value = 216635.76
$ws.Range("M141").Value = -211455.76

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 150836.64
$ws.Range("I32").Value = 144545.81
$ws.Range("K32").Value = 144545.81
$ws.Range("M32").Value = -144258.81

$ws.Range("H63").Value = 1237.5
$ws.Range("I63").Value = 1185
$ws.Range("K63").Value = 1185
$ws.Range("M63").Value = -499

$ws.Range("H66").Value = 1237.5
$ws.Range("I66").Value = 1185
$ws.Range("K66").Value = 5925
$ws.Range("M66").Value = -2493

$ws.Range("H132").Value = 1993.3
$ws.Range("I132").Value = 2014.7778
$ws.Range("K132").Value = 6044.3334
$ws.Range("M132").Value = -3514.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H110").Value = 78455.75
$ws.Range("J110").Value = 78455.75
$ws.Range("L110").Value = 78455.75
$ws.Range("N110").Value = -86635.75

$ws.Range("H134").Value = 2260.6086
$ws.Range("I134").Value = 1840.3158
$ws.Range("J134").Value = 4257
$ws.Range("K134").Value = 5520.9474
$ws.Range("L134").Value = 12771
$ws.Range("M134").Value = -2985.9474
$ws.Range("N134").Value = -17841

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1297.25
$ws.Range("I68").Value = 1070.1428
$ws.Range("J68").Value = 1615.2
$ws.Range("K68").Value = 3210.4284
$ws.Range("L68").Value = 4845.6
$ws.Range("M68").Value = -2399.4284
$ws.Range("N68").Value = -6467.6

$ws.Range("H71").Value = 1297.25
$ws.Range("I71").Value = 1070.1428
$ws.Range("J71").Value = 1615.2
$ws.Range("K71").Value = 9631.2852
$ws.Range("L71").Value = 14536.8
$ws.Range("M71").Value = -5575.2852
$ws.Range("N71").Value = -22648.8

$ws.Range("H80").Value = 2165.6667
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()

$ws.Range("H83").Value = 2165.6667
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()

$ws.Range("H86").Value = 396
$ws.Range("I86").Value = 396
$ws.Range("K86").Value = 1188
$ws.Range("M86").Value = -2

$ws.Range("H89").Value = 396
$ws.Range("I89").Value = 396
$ws.Range("K89").Value = 3564
$ws.Range("M89").Value = 2364

$ws.Range("H97").Value = 1746.7
$ws.Range("J97").Value = 1939.625
$ws.Range("L97").Value = 5818.875
$ws.Range("N97").Value = -6810.875

$ws.Range("H98").Value = 416
$ws.Range("I98").Value = 254.16667
$ws.Range("J98").Value = 610.2
$ws.Range("K98").Value = 762.50001
$ws.Range("L98").Value = 1830.6
$ws.Range("M98").Value = 735.49999
$ws.Range("N98").Value = -4826.6

$ws.Range("H113").Value = 641.4286
$ws.Range("I113").Value = 578
$ws.Range("K113").Value = 1734
$ws.Range("M113").Value = 436

$ws.Range("H114").Value = 20001272
$ws.Range("I114").Value = 40001044
$ws.Range("J114").Value = 1500
$ws.Range("K114").Value = 120003132
$ws.Range("L114").Value = 4500
$ws.Range("M114").Value = -119999878
$ws.Range("N114").Value = -11008

$ws.Range("H121").Value = 33335818
$ws.Range("J121").Value = 4326.2
$ws.Range("L121").Value = 12978.6
$ws.Range("N121").Value = -15598.6

$ws.Range("H129").Value = 2214.0557
$ws.Range("J129").Value = 2214.0557
$ws.Range("L129").Value = 6642.1671
$ws.Range("N129").Value = -16642.1671

$ws.Range("H130").Value = 2999
$ws.Range("I130").Value = 2999
$ws.Range("J130").Value = 0
$ws.Range("K130").Value = 8997
$ws.Range("L130").Value = 0
$ws.Range("M130").Value = -3977
$ws.Range("N130").ClearContents()

$ws.Range("H131").Value = 2224.6667
$ws.Range("I131").Value = 1684.8
$ws.Range("J131").Value = 2494.6
$ws.Range("K131").Value = 5054.4
$ws.Range("L131").Value = 7483.799999999999
$ws.Range("M131").Value = -14.39999999999964
$ws.Range("N131").Value = -17563.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 27804444
$ws.Range("I24").Value = 83346664
$ws.Range("J24").Value = 33333
$ws.Range("K24").Value = 83346664
$ws.Range("L24").Value = 33333
$ws.Range("M24").Value = -83346491
$ws.Range("N24").Value = -33679

$ws.Range("H70").Value = 9195.1
$ws.Range("I70").Value = 8281.714
$ws.Range("K70").Value = 8281.714
$ws.Range("M70").Value = -8011.714

$ws.Range("H73").Value = 9195.1
$ws.Range("I73").Value = 8281.714
$ws.Range("K73").Value = 8281.714
$ws.Range("M73").Value = -7345.714

$ws.Range("H97").Value = 19445.371
$ws.Range("I97").Value = 26830.44
$ws.Range("J97").Value = 982.7
$ws.Range("K97").Value = 26830.44
$ws.Range("L97").Value = 982.7
$ws.Range("M97").Value = -26334.44
$ws.Range("N97").Value = -1974.7

$ws.Range("H132").Value = 1860.5652
$ws.Range("I132").Value = 1746.9445
$ws.Range("J132").Value = 2269.6
$ws.Range("K132").Value = 5240.833500000001
$ws.Range("L132").Value = 6808.799999999999
$ws.Range("M132").Value = -2710.833500000001
$ws.Range("N132").Value = -11868.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 40000
$ws.Range("J4").Value = 40000
$ws.Range("L4").Value = 40000
$ws.Range("N4").Value = -40226

$ws.Range("H28").Value = 40000
$ws.Range("J28").Value = 40000
$ws.Range("L28").Value = 40000
$ws.Range("N28").Value = -40464

$ws.Range("H37").Value = 40000
$ws.Range("J37").Value = 40000
$ws.Range("L37").Value = 40000
$ws.Range("N37").Value = -40214

$ws.Range("H55").Value = 242.91667
$ws.Range("I55").Value = 143.54546
$ws.Range("J55").Value = 327
$ws.Range("K55").Value = 143.54546
$ws.Range("L55").Value = 327
$ws.Range("M55").Value = 29.45454000000001
$ws.Range("N55").Value = -673

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("J4").Value = 40000
$ws.Range("L4").Value = 40000
$ws.Range("N4").Value = -40226

$ws.Range("H21").Value = 2500
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 2500
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 2500
$ws.Range("M21").ClearContents()
$ws.Range("N21").Value = -2970

$ws.Range("H26").Value = 12500
$ws.Range("J26").Value = 12500
$ws.Range("L26").Value = 12500
$ws.Range("N26").Value = -13086

$ws.Range("H35").Value = 2500
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 2500
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 2500
$ws.Range("M35").ClearContents()
$ws.Range("N35").Value = -3080

$ws.Range("H113").Value = 918.7143
$ws.Range("I113").Value = 857.4
$ws.Range("J113").Value = 1072
$ws.Range("K113").Value = 2572.2
$ws.Range("L113").Value = 3216
$ws.Range("M113").Value = -402.1999999999998
$ws.Range("N113").Value = -7556
